# "added some checks and consistency to the Excel file"
#
# - C3 (Time(s) for the first date row) is corrected from 1 to 20.
# - The extra sample rows (4-11) are no longer needed: rows 4 and 5 are
#   cleared but kept, rows 6-11's old Date/Time(s) pairs are removed, which
#   shrinks the sheet's used range down to row 6.
# - A new, empty "marker" cell E6 is added with an underlined font, as a
#   small consistency/formatting check mark near the trimmed data.
# - The A2:C3 header+first-row block is left selected for review.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Time(s) value for the first data row.
$ws.Range("C3").Value = 20

# Drop the now-unneeded sample rows below it (was rows 4-11, Date/Time(s)
# pairs 2 through 9).
$ws.Range("B4:C11").ClearContents()

# Add a styled (underlined) empty check-cell a couple of columns over.
$ws.Range("E6").Font.Underline = $true

# Leave the header/first-row block selected.
$ws.Range("A2:C3").Select()
